$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value = 50052
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 50052
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 50052
$ws.Range("N126").Value = -59932

$ws.Range("H132").Value = 1631.25
$ws.Range("I132").Value = 1714.6666
$ws.Range("J132").Value = 380
$ws.Range("K132").Value = 5143.9998
$ws.Range("L132").Value = 1140
$ws.Range("M132").Value = -2613.9998
$ws.Range("N132").Value = -6200

$ws.Range("H138").Value = 1591121.6
$ws.Range("I138").Value = 1803.6538
$ws.Range("J138").Value = 2707939.5
$ws.Range("K138").Value = 5410.9614
$ws.Range("L138").Value = 8123818.5
$ws.Range("M138").Value = -270.9614000000001
$ws.Range("N138").Value = -8134098.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2344.2856
$ws.Range("I2").Value = 1710.9166
$ws.Range("J2").Value = 3188.7778
$ws.Range("K2").Value = 1710.9166
$ws.Range("L2").Value = 3188.7778
$ws.Range("M2").Value = -1597.9166
$ws.Range("N2").Value = -3414.7778

$ws.Range("H61").Value = 8562.526
$ws.Range("I61").Value = 2196.5557
$ws.Range("J61").Value = 14291.9
$ws.Range("K61").Value = 2196.5557
$ws.Range("L61").Value = 14291.9
$ws.Range("M61").Value = -1984.5557
$ws.Range("N61").Value = -14715.9

$ws.Range("H116").Value = 2344.2856
$ws.Range("I116").Value = 1710.9166
$ws.Range("J116").Value = 3188.7778
$ws.Range("K116").Value = 1710.9166
$ws.Range("L116").Value = 3188.7778
$ws.Range("M116").Value = 583.0834
$ws.Range("N116").Value = -7776.7778

$ws.Range("H136").Value = 8562.526
$ws.Range("I136").Value = 2196.5557
$ws.Range("J136").Value = 14291.9
$ws.Range("K136").Value = 6589.6671
$ws.Range("L136").Value = 42875.7
$ws.Range("M136").Value = -4039.6671
$ws.Range("N136").Value = -47975.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2344.2856
$ws.Range("I3").Value = 1710.9166
$ws.Range("J3").Value = 3188.7778
$ws.Range("K3").Value = 1710.9166
$ws.Range("L3").Value = 3188.7778
$ws.Range("M3").Value = -1596.9166
$ws.Range("N3").Value = -3416.7778

$ws.Range("H26").Value = 33601
$ws.Range("I26").Value = 26820.666
$ws.Range("J26").Value = 53942
$ws.Range("K26").Value = 26820.666
$ws.Range("L26").Value = 53942
$ws.Range("M26").Value = -26528.666
$ws.Range("N26").Value = -54526

$ws.Range("H62").Value = 42999.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 42999.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 42999.5
$ws.Range("N62").Value = -44371.5

$ws.Range("H65").Value = 42999.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 42999.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 128998.5
$ws.Range("N65").Value = -135862.5

$ws.Range("H107").Value = 43274120
$ws.Range("I107").Value = 46880060
$ws.Range("J107").Value = 2845
$ws.Range("K107").Value = 46880060
$ws.Range("L107").Value = 2845
$ws.Range("M107").Value = -46878140
$ws.Range("N107").Value = -6685

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6131.618
$ws.Range("I31").Value = 1659.96
$ws.Range("J31").Value = 9858
$ws.Range("K31").Value = 1659.96
$ws.Range("L31").Value = 9858
$ws.Range("M31").Value = -1364.96
$ws.Range("N31").Value = -10448

$ws.Range("H34").Value = 6131.618
$ws.Range("I34").Value = 1659.96
$ws.Range("J34").Value = 9858
$ws.Range("K34").Value = 1659.96
$ws.Range("L34").Value = 9858
$ws.Range("M34").Value = -1457.96
$ws.Range("N34").Value = -10262

$ws.Range("H58").Value = 5727.4224
$ws.Range("I58").Value = 2038.3334
$ws.Range("J58").Value = 8186.815
$ws.Range("K58").Value = 2038.3334
$ws.Range("L58").Value = 8186.815
$ws.Range("M58").Value = -1835.3334
$ws.Range("N58").Value = -8592.814999999999

$ws.Range("H134").Value = 5822.6333
$ws.Range("I134").Value = 1622.75
$ws.Range("J134").Value = 10622.5
$ws.Range("K134").Value = 4868.25
$ws.Range("L134").Value = 31867.5
$ws.Range("M134").Value = -2333.25
$ws.Range("N134").Value = -36937.5

$ws.Range("H136").Value = 5727.4224
$ws.Range("I136").Value = 2038.3334
$ws.Range("J136").Value = 8186.815
$ws.Range("K136").Value = 6115.0002
$ws.Range("L136").Value = 24560.445
$ws.Range("M136").Value = -3565.0002
$ws.Range("N136").Value = -29660.445

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4208.524
$ws.Range("I126").Value = 2428.4167
$ws.Range("J126").Value = 6582
$ws.Range("K126").Value = 7285.250100000001
$ws.Range("L126").Value = 19746
$ws.Range("M126").Value = -4815.250100000001
$ws.Range("N126").Value = -24686

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1590.2142
$ws.Range("I22").Value = 487.8
$ws.Range("J22").Value = 4346.25
$ws.Range("K22").Value = 487.8
$ws.Range("L22").Value = 4346.25
$ws.Range("M22").Value = -192.8
$ws.Range("N22").Value = -4936.25

$ws.Range("H27").Value = 1590.2142
$ws.Range("I27").Value = 487.8
$ws.Range("J27").Value = 4346.25
$ws.Range("K27").Value = 487.8
$ws.Range("L27").Value = 4346.25
$ws.Range("M27").Value = -380.8
$ws.Range("N27").Value = -4560.25

$ws.Range("H46").Value = 2375.9348
$ws.Range("I46").Value = 1958.1389
$ws.Range("J46").Value = 3880
$ws.Range("K46").Value = 1958.1389
$ws.Range("L46").Value = 3880
$ws.Range("M46").Value = -1770.1389
$ws.Range("N46").Value = -4256

$ws.Range("H55").Value = 604
$ws.Range("I55").Value = 160.875
$ws.Range("J55").Value = 876.6923
$ws.Range("K55").Value = 160.875
$ws.Range("L55").Value = 876.6923
$ws.Range("M55").Value = 12.125
$ws.Range("N55").Value = -1222.6923

$ws.Range("H61").Value = 4537.5747
$ws.Range("I61").Value = 3278.0938
$ws.Range("J61").Value = 7224.467
$ws.Range("K61").Value = 3278.0938
$ws.Range("L61").Value = 7224.467
$ws.Range("M61").Value = -3076.0938
$ws.Range("N61").Value = -7628.467

$ws.Range("H68").Value = 6564.4287
$ws.Range("I68").Value = 7750
$ws.Range("J68").Value = 6366.8335
$ws.Range("K68").Value = 7750
$ws.Range("L68").Value = 6366.8335
$ws.Range("M68").Value = -7001
$ws.Range("N68").Value = -7864.8335

$ws.Range("H71").Value = 6564.4287
$ws.Range("I71").Value = 7750
$ws.Range("J71").Value = 6366.8335
$ws.Range("K71").Value = 38750
$ws.Range("L71").Value = 31834.1675
$ws.Range("M71").Value = -35006
$ws.Range("N71").Value = -39322.1675

$ws.Range("H82").Value = 56112420
$ws.Range("I82").Value = 84167256
$ws.Range("J82").Value = 2747.1667
$ws.Range("K82").Value = 84167256
$ws.Range("L82").Value = 2747.1667
$ws.Range("M82").Value = -84166895
$ws.Range("N82").Value = -3469.1667

$ws.Range("H85").Value = 56112420
$ws.Range("I85").Value = 84167256
$ws.Range("J85").Value = 2747.1667
$ws.Range("K85").Value = 84167256
$ws.Range("L85").Value = 2747.1667
$ws.Range("M85").Value = -84166008
$ws.Range("N85").Value = -5243.1667

$ws.Range("H100").Value = 4264.1
$ws.Range("I100").Value = 3270.0833
$ws.Range("J100").Value = 5755.125
$ws.Range("K100").Value = 3270.0833
$ws.Range("L100").Value = 5755.125
$ws.Range("M100").Value = -2729.0833
$ws.Range("N100").Value = -6837.125

$ws.Range("H113").Value = 4537.5747
$ws.Range("I113").Value = 3278.0938
$ws.Range("J113").Value = 7224.467
$ws.Range("K113").Value = 3278.0938
$ws.Range("L113").Value = 7224.467
$ws.Range("M113").Value = -1108.0938
$ws.Range("N113").Value = -11564.467

$ws.Range("H136").Value = 12104.1455
$ws.Range("I136").Value = 3025
$ws.Range("J136").Value = 18589.25
$ws.Range("K136").Value = 9075
$ws.Range("L136").Value = 55767.75
$ws.Range("M136").Value = -6525
$ws.Range("N136").Value = -60867.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 133339200
$ws.Range("I62").Value = 142862380
$ws.Range("J62").Value = 111118450
$ws.Range("K62").Value = 142862380
$ws.Range("L62").Value = 111118450
$ws.Range("M62").Value = -142861756
$ws.Range("N62").Value = -111119698

$ws.Range("H65").Value = 133339200
$ws.Range("I65").Value = 142862380
$ws.Range("J65").Value = 111118450
$ws.Range("K65").Value = 714311900
$ws.Range("L65").Value = 555592250
$ws.Range("M65").Value = -714308780
$ws.Range("N65").Value = -555598490

$ws.Range("H126").Value = 1193.5385
$ws.Range("I126").Value = 1190.7142
$ws.Range("J126").Value = 1196.8334
$ws.Range("K126").Value = 3572.1426
$ws.Range("L126").Value = 3590.5002
$ws.Range("M126").Value = -1102.1426
$ws.Range("N126").Value = -8530.5002

$ws.Range("H132").Value = 2983.9546
$ws.Range("I132").Value = 1260.2142
$ws.Range("J132").Value = 6000.5
$ws.Range("K132").Value = 3780.6426
$ws.Range("L132").Value = 18001.5
$ws.Range("M132").Value = -1250.6426
$ws.Range("N132").Value = -23061.5

$ws.Range("H136").Value = 4485.1943
$ws.Range("I136").Value = 1659.0667
$ws.Range("J136").Value = 6503.857
$ws.Range("K136").Value = 4977.2001
$ws.Range("L136").Value = 19511.571
$ws.Range("M136").Value = -2427.2001
$ws.Range("N136").Value = -24611.571
